# Add two new slides to the end of the deck: "Reproducible Report
# Location" and "Thank You!" - both built on the "Title and Content"
# layout (CustomLayout index 2), the same layout already used by the
# other content slides in this deck (e.g. slide 23, "Sources of Bias").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 24: "Reproducible Report Location"
# ---------------------------------------------------------------------
$slide24 = $p.Slides.Add($p.Slides.Count + 1, 2)

$title24 = $slide24.Shapes.Placeholders.Item(1)
$title24.TextFrame.TextRange.Text = "Reproducible Report Location"

$body24 = $slide24.Shapes.Placeholders.Item(2)

# Reposition / resize the content placeholder to match the authored slide.
$body24.Left = 1103312 / 12700
$body24.Top = 2052918 / 12700
$body24.Width = 10629652 / 12700
$body24.Height = 4195481 / 12700

$tr24 = $body24.TextFrame.TextRange
$tr24.Text = "https://"
[void]$tr24.InsertAfter("github.com/clayv/UC_Boulder_DTSA5301_Final/tree/main/NYPD_Shootings")
[void]$tr24.InsertAfter("`rContains")
[void]$tr24.InsertAfter("`rRStudio")
[void]$tr24.InsertAfter(" ")
[void]$tr24.InsertAfter("Rmd")
[void]$tr24.InsertAfter("`rKnitted HTML")
[void]$tr24.InsertAfter("`rPowerpoint")
[void]$tr24.InsertAfter(" slide deck")

# Indent the three "Contains" sub-bullets to the second outline level.
$tr24.Paragraphs(3).IndentLevel = 2
$tr24.Paragraphs(4).IndentLevel = 2
$tr24.Paragraphs(5).IndentLevel = 2

# ---------------------------------------------------------------------
# Slide 25: "Thank You!"
# ---------------------------------------------------------------------
$slide25 = $p.Slides.Add($p.Slides.Count + 1, 2)

$title25 = $slide25.Shapes.Placeholders.Item(1)
$title25.TextFrame.TextRange.Text = "Thank You!"

# The authored slide has no content placeholder on it - remove the
# empty one that comes by default with the layout.
$slide25.Shapes.Placeholders.Item(2).Delete()
